$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows and columns (letters) where a new "A" (absent) mark must be entered.
# Each of these cells goes from an untouched/empty "s=38" style cell to a
# "s=39" styled cell holding the shared string "A" (index 20).
$marks = @(
    "L10",
    "K11",
    "K14",
    "K15",
    "K17",
    "K19", "L19",
    "K25",
    "K26",
    "K27", "L27",
    "K28",
    "K29",
    "L33",
    "L34",
    "L42",
    "K46",
    "K53",
    "K54",
    "K57",
    "L72",
    "L76",
    "K83"
)

foreach ($addr in $marks) {
    $c = $ws.Range($addr)
    $c.Value = "A"
    # Touching an alignment property (leaving it at its default) bumps the
    # cell onto the "touched" style variant (s=39) without altering the
    # visible formatting, matching what Excel does when a cell is edited.
    $c.WrapText = $false
}

# J83 and L20 receive the same "touched" style bump but their content stays
# blank (the attendance mark was toggled on and back off again).
$ws.Range("J83").WrapText = $false
$ws.Range("L20").WrapText = $false

# Change the frozen panes: freeze top 3 rows + first 3 columns, active cell D4.
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("D4").Select()
$excel.ActiveWindow.FreezePanes = $true
